$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value2 = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextCell "D2" '52.214.11'
Set-TextCell "E2" '  +0.93%  '

# Row 3
Set-TextCell "D3" '2.910.82'
Set-TextCell "E3" '  +3.92%  '

# Row 5
Set-TextCell "D5" '351.46'
Set-TextCell "E5" '  -1.31%  '

# Row 6
Set-TextCell "D6" '112.35'
Set-TextCell "E6" '  +3.01%  '

# Row 7
Set-TextCell "D7" '0.558'
Set-TextCell "E7" '  +0.86%  '

# Row 8
Set-TextCell "E8" '  +0.12%  '

# Row 9
Set-TextCell "D9" '0.633'
Set-TextCell "E9" '  +0.51%  '

# Row 10
Set-TextCell "D10" '39.92'
Set-TextCell "E10" '  +0.12%  '

# Row 11
Set-TextCell "D11" '0.0861'
Set-TextCell "E11" '  +2.59%  '

# Row 12
Set-TextCell "E12" '  +0.36%  '

# Row 13
Set-TextCell "D13" '20.00'
Set-TextCell "E13" '  +0.51%  '

# Row 14
Set-TextCell "D14" '7.80'
Set-TextCell "E14" '  +0.46%  '

# Row 15
Set-TextCell "D15" '3.371.52'
Set-TextCell "E15" '  +3.95%  '

# Row 16
Set-TextCell "B16" 'Polygon'
Set-TextCell "C16" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextCell "D16" '0.998'
Set-TextCell "E16" '  +5.97%  '

# Row 17
Set-TextCell "B17" 'WrappedEther'
Set-TextCell "C17" 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextCell "D17" '2.899.17'
Set-TextCell "E17" '  +3.35%  '

# Row 18
Set-TextCell "D18" '52.266.72'
Set-TextCell "E18" '  +1.11%  '

# Row 19
Set-TextCell "D19" '7.64'
Set-TextCell "E19" '  -0.57%  '

# Row 20
Set-TextCell "E20" '  +5.08%  '

# Row 21
Set-TextCell "D21" '14.18'
Set-TextCell "E21" '  +4.54%  '

# Row 22
Set-TextCell "D22" '0.0₃0980'
Set-TextCell "E22" '  +0.24%  '

# Row 23
Set-TextCell "D23" '70.89'
Set-TextCell "E23" '  +0.74%  '

# Row 24
Set-TextCell "D24" '269.75'
Set-TextCell "E24" '  +0.56%  '

# Row 25
Set-TextCell "E25" '  +1.93%  '

# Row 26
Set-TextCell "D26" '26.73'
Set-TextCell "E26" '  +2.39%  '

# Row 27
Set-TextCell "D27" '0.999'
Set-TextCell "E27" '  -0.08%  '

# Row 28
Set-TextCell "D28" '0.164'
Set-TextCell "E28" '  +0.37%  '

# Row 29
Set-TextCell "D29" '10.59'
Set-TextCell "E29" '  +2.20%  '

# Row 30
Set-TextCell "D30" '37.62'
Set-TextCell "E30" '  -0.04%  '

# Row 31
Set-TextCell "D31" '6.50'
Set-TextCell "E31" '  +4.50%  '

# Row 32
Set-TextCell "D32" '2.25'
Set-TextCell "E32" '  +0.91%  '

# Row 33
Set-TextCell "D33" '6.14'
Set-TextCell "E33" '  +7.84%  '

# Row 34
Set-TextCell "D34" '0.0955'
Set-TextCell "E34" '  +11.07%  '

# Row 35
Set-TextCell "D35" '53.26'
Set-TextCell "E35" '  +2.56%  '

# Row 36
Set-TextCell "E36" '  +2.09%  '

# Row 37
Set-TextCell "D37" '0.999'
Set-TextCell "E37" '  -0.13%  '

# Row 38
Set-TextCell "D38" '3.30'
Set-TextCell "E38" '  +5.23%  '

# Row 39
Set-TextCell "D39" '2.07'
Set-TextCell "E39" '  +3.71%  '

# Row 40
Set-TextCell "D40" '18.69'
Set-TextCell "E40" '  -0.63%  '

# Row 41
Set-TextCell "D41" '2.82'
Set-TextCell "E41" '  +13.51%  '

# Row 42
Set-TextCell "E42" '  +1.40%  '

# Row 43
Set-TextCell "E43" '  +6.99%  '

# Row 44
Set-TextCell "B44" 'Monero'
Set-TextCell "C44" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell "D44" '121.33'
Set-TextCell "E44" '  +1.88%  '

# Row 45
Set-TextCell "B45" 'ApeXProtocol'
Set-TextCell "C45" 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell "D45" '2.63'
Set-TextCell "E45" '  +7.36%  '

# Row 46
Set-TextCell "E46" '  -0.93%  '

# Row 47
Set-TextCell "D47" '2.199.87'
Set-TextCell "E47" '  +4.42%  '

# Row 48
Set-TextCell "D48" '3.52'
Set-TextCell "E48" '  +4.22%  '

# Row 49
Set-TextCell "D49" '0.261'
Set-TextCell "E49" '  +23.04%  '

# Row 50
Set-TextCell "E50" '  +12.02%  '

# Row 51
Set-TextCell "E51" '  +5.81%  '
